# This workbook keeps a rolling weekly price history. A new weekly entry is
# inserted at row 23 (pushing the existing historical rows 23-282 down to
# 24-283), and the rest of the row's data (Volumen/Precio columns, etc.) is
# carried over from the row that used to be in that slot - only the date
# (column D) is updated to the new reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23; Excel automatically shifts rows
# 23:282 down to 24:283 and updates the sheet dimension accordingly.
$ws.Rows("23:23").Insert()

# The newly inserted row 23 is empty. Populate it with the same data that is
# now sitting in row 24 (i.e. what used to be the old row 23), then overwrite
# just the date so the new row reflects the new weekly reading.
$ws.Range("A24:R24").Copy($ws.Range("A23:R23"))
$ws.Range("D23").Value = 44630
